$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 34966.668
$ws.Range("J3").Value = 34966.668
$ws.Range("L3").Value = 34966.668
$ws.Range("N3").Value = -35194.668

$ws.Range("H18").Value = 1245
$ws.Range("I18").Value = 990
$ws.Range("K18").Value = 990
$ws.Range("M18").Value = -706

$ws.Range("H102").Value = 34966.668
$ws.Range("J102").Value = 34966.668
$ws.Range("L102").Value = 34966.668
$ws.Range("N102").Value = -41456.668

$ws.Range("I127").Value = 3799.4
$ws.Range("J127").Value = 1916.5
$ws.Range("K127").Value = 11398.2
$ws.Range("L127").Value = 5749.5
$ws.Range("M127").Value = -6438.200000000001
$ws.Range("N127").Value = -15669.5

$ws.Range("H137").Value = 1616
$ws.Range("J137").Value = 2884.5
$ws.Range("L137").Value = 8653.5
$ws.Range("N137").Value = -13753.5

$ws.Range("H141").Value = 416.33334
$ws.Range("I141").Value = 416.33334
$ws.Range("K141").Value = 1249.00002
$ws.Range("M141").Value = 3930.99998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 750
$ws.Range("I45").Value = 750
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 750
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -373
$ws.Range("N45").ClearContents()

$ws.Range("H103").Value = 32000
$ws.Range("J103").Value = 32000
$ws.Range("L103").Value = 32000
$ws.Range("N103").Value = -34344

$ws.Range("H110").Value = 1238
$ws.Range("I110").Value = 1356.5714
$ws.Range("K110").Value = 1356.5714
$ws.Range("M110").Value = 688.4286

$ws.Range("H132").Value = 790.5
$ws.Range("I132").Value = 790.5
$ws.Range("K132").Value = 2371.5
$ws.Range("M132").Value = 158.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H31").Value = 25011.5
$ws.Range("I31").Value = 9023
$ws.Range("J31").Value = 41000
$ws.Range("K31").Value = 9023
$ws.Range("L31").Value = 41000
$ws.Range("M31").Value = -8771
$ws.Range("N31").Value = -41504

$ws.Range("H81").Value = 39822.5
$ws.Range("J81").Value = 39822.5
$ws.Range("L81").Value = 39822.5
$ws.Range("N81").Value = -41944.5

$ws.Range("H84").Value = 39822.5
$ws.Range("J84").Value = 39822.5
$ws.Range("L84").Value = 119467.5
$ws.Range("N84").Value = -130075.5

$ws.Range("H88").Value = 21562.25
$ws.Range("J88").Value = 21562.25
$ws.Range("L88").Value = 21562.25
$ws.Range("N88").Value = -22374.25

$ws.Range("H91").Value = 21562.25
$ws.Range("J91").Value = 21562.25
$ws.Range("L91").Value = 21562.25
$ws.Range("N91").Value = -24370.25

$ws.Range("H103").Value = 12000
$ws.Range("J103").Value = 12000
$ws.Range("L103").Value = 12000
$ws.Range("N103").Value = -14344

$ws.Range("H106").Value = 5677
$ws.Range("J106").Value = 5677
$ws.Range("L106").Value = 5677
$ws.Range("N106").Value = -8201

$ws.Range("H123").Value = 100780
$ws.Range("J123").Value = 100780
$ws.Range("L123").Value = 100780
$ws.Range("N123").Value = -110580

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 9000
$ws.Range("I16").Value = 9000
$ws.Range("K16").Value = 9000
$ws.Range("M16").Value = -8713

$ws.Range("H31").Value = 1434.8948
$ws.Range("I31").Value = 1441.4445
$ws.Range("J31").Value = 1429
$ws.Range("K31").Value = 1441.4445
$ws.Range("L31").Value = 1429
$ws.Range("M31").Value = -1146.4445
$ws.Range("N31").Value = -2019

$ws.Range("H34").Value = 1434.8948
$ws.Range("I34").Value = 1441.4445
$ws.Range("J34").Value = 1429
$ws.Range("K34").Value = 1441.4445
$ws.Range("L34").Value = 1429
$ws.Range("M34").Value = -1239.4445
$ws.Range("N34").Value = -1833

$ws.Range("H43").Value = 12749.75
$ws.Range("J43").Value = 12749.75
$ws.Range("L43").Value = 12749.75
$ws.Range("N43").Value = -13117.75

$ws.Range("H101").Value = 12749.75
$ws.Range("J101").Value = 12749.75
$ws.Range("L101").Value = 12749.75
$ws.Range("N101").Value = -19239.75

$ws.Range("H113").Value = 9000
$ws.Range("I113").Value = 9000
$ws.Range("K113").Value = 9000
$ws.Range("M113").Value = -6830

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 18.777779
$ws.Range("I2").Value = 17.375
$ws.Range("K2").Value = 104.25
$ws.Range("M2").Value = 8.75

$ws.Range("H6").Value = 111151.336
$ws.Range("I6").Value = 47
$ws.Range("J6").Value = 250031.75
$ws.Range("K6").Value = 141
$ws.Range("L6").Value = 750095.25
$ws.Range("M6").Value = -28
$ws.Range("N6").Value = -750321.25

$ws.Range("H7").Value = 209
$ws.Range("I7").Value = 357
$ws.Range("J7").Value = 98
$ws.Range("K7").Value = 1071
$ws.Range("L7").Value = 294
$ws.Range("M7").Value = -959
$ws.Range("N7").Value = -518

$ws.Range("H23").Value = 520.7143
$ws.Range("I23").Value = 374
$ws.Range("J23").Value = 545.1667
$ws.Range("K23").Value = 1122
$ws.Range("L23").Value = 1635.5001
$ws.Range("M23").Value = -887
$ws.Range("N23").Value = -2105.5001

$ws.Range("H34").Value = 1744.2632
$ws.Range("I34").Value = 95.46154
$ws.Range("J34").Value = 5316.6665
$ws.Range("K34").Value = 286.38462
$ws.Range("L34").Value = 15949.9995
$ws.Range("M34").Value = -202.38462
$ws.Range("N34").Value = -16117.9995

$ws.Range("H98").Value = 275.84616
$ws.Range("I98").Value = 347.5
$ws.Range("J98").Value = 161.2
$ws.Range("K98").Value = 1042.5
$ws.Range("L98").Value = 483.6
$ws.Range("M98").Value = 455.5
$ws.Range("N98").Value = -3479.6

$ws.Range("H107").Value = 1100
$ws.Range("I107").Value = 1000
$ws.Range("K107").Value = 3000
$ws.Range("M107").Value = -1080

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 100
$ws.Range("I27").Value = 100
$ws.Range("J27").Value = 100
$ws.Range("K27").Value = 100
$ws.Range("L27").Value = 100
$ws.Range("M27").Value = 66
$ws.Range("N27").Value = -432

$ws.Range("H41").Value = 1025.5
$ws.Range("I41").Value = 1025.5
$ws.Range("K41").Value = 1025.5
$ws.Range("M41").Value = -670.5

$ws.Range("H101").Value = 19998
$ws.Range("J101").Value = 19998
$ws.Range("L101").Value = 19998
$ws.Range("N101").Value = -26488

$ws.Range("H104").Value = 40000
$ws.Range("J104").Value = 40000
$ws.Range("L104").Value = 40000
$ws.Range("N104").Value = -46988

$ws.Range("H105").Value = 15999
$ws.Range("J105").Value = 15999
$ws.Range("L105").Value = 15999
$ws.Range("N105").Value = -22987

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 420.5
$ws.Range("I22").Value = 420.5
$ws.Range("K22").Value = 420.5
$ws.Range("M22").Value = -125.5

$ws.Range("H27").Value = 420.5
$ws.Range("I27").Value = 420.5
$ws.Range("K27").Value = 420.5
$ws.Range("M27").Value = -313.5

$ws.Range("H30").Value = 14284.5
$ws.Range("I30").Value = 560
$ws.Range("J30").Value = 28009
$ws.Range("K30").Value = 560
$ws.Range("L30").Value = 28009
$ws.Range("M30").Value = -452
$ws.Range("N30").Value = -28225

$ws.Range("H35").Value = 2177.3333
$ws.Range("I35").Value = 1765.5
$ws.Range("K35").Value = 1765.5
$ws.Range("M35").Value = -1429.5

$ws.Range("H39").Value = 59
$ws.Range("I39").Value = 59
$ws.Range("K39").Value = 59
$ws.Range("M39").Value = 401

$ws.Range("H58").Value = 9000
$ws.Range("I58").Value = 8000
$ws.Range("K58").Value = 8000
$ws.Range("M58").Value = -7740

$ws.Range("H132").Value = 2000
$ws.Range("I132").Value = 2000
$ws.Range("K132").Value = 6000
$ws.Range("M132").Value = -3470

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()

$ws.Range("H98").Value = 13000
$ws.Range("J98").Value = 13000
$ws.Range("L98").Value = 13000
$ws.Range("N98").Value = -18990

$ws.Range("H107").Value = 1091.3334
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 1091.3334
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 3274.0002
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -7114.0002

$ws.Range("H136").Value = 2416
$ws.Range("I136").Value = 1124.25
$ws.Range("K136").Value = 3372.75
$ws.Range("M136").Value = -822.75
